# Change Minedu name in school transports templates body
#
# The ministry name "Παιδείας, Έρευνας και Θρησκευμάτων"
# (Ministry of Education, Research and Religious Affairs) is updated to
# "Παιδείας και Θρησκευμάτων" (Ministry of Education and Religious Affairs)
# in the body paragraph describing the post-mobility reporting obligation.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Υπουργείο  Παιδείας, Έρευνας και Θρησκευμάτων στην Αυτοτελή",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Υπουργείο  Παιδείας και Θρησκευμάτων στην Αυτοτελή",
    2
)
